$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear columns AB-AK and AM for rows 2 through 34, leaving AL untouched.
$ws.Range("AB2:AK34").ClearContents()
$ws.Range("AM2:AM34").ClearContents()
